{"js": "// Add null-support block for the multi-line \"@if\" / \"@endif\" sample:\n// duplicate the last @if ... / ... / @endif block, but switch the\n// condition to \"CreateDate == NULL\", after a blank separator paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The very last paragraph in the document is an empty trailing paragraph;\n// insert the new block right before it, directly after the existing\n// \"@endif\" paragraph.\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\nlet anchor = lastParagraph.insertParagraph(\"\", Word.InsertLocation.before);\nanchor = anchor.insertParagraph(\"@if CreateDate == NULL\", Word.InsertLocation.after);\nanchor = anchor.insertParagraph(\"CreateDate is not less than 2021\", Word.InsertLocation.after);\nanchor = anchor.insertParagraph(\"asdasdasdasdasdasd\", Word.InsertLocation.after);\nanchor = anchor.insertParagraph(\"@endif\", Word.InsertLocation.after);\n\nawait context.sync();\n", "ps1": "# Add null-support block for the multi-line \"@if\" / \"@endif\" sample:\n# duplicate the last @if ... / ... / @endif block, but switch the\n# condition to \"CreateDate == NULL\", after a blank separator paragraph.\n\n$d = $word.ActiveDocument\n\n$count = $d.Paragraphs.Count\n$endifPara = $d.Paragraphs.Item($count - 1)   # existing \"@endif\" paragraph\n\n# Blank separator paragraph.\n$endifPara.Range.InsertParagraphAfter()\n$blank = $d.Paragraphs.Item($count)\n\n# \"@if CreateDate == NULL\"\n$blank.Range.InsertParagraphAfter()\n$pIf = $d.Paragraphs.Item($count + 1)\n$pIf.Range.InsertAfter(\"@if CreateDate == NULL\")\n\n# \"CreateDate is not less than 2021\"\n$pIf.Range.InsertParagraphAfter()\n$pBody = $d.Paragraphs.Item($count + 2)\n$pBody.Range.InsertAfter(\"CreateDate is not less than 2021\")\n\n# \"asdasdasdasdasdasd\"\n$pBody.Range.InsertParagraphAfter()\n$pFiller = $d.Paragraphs.Item($count + 3)\n$pFiller.Range.InsertAfter(\"asdasdasdasdasdasd\")\n\n# \"@endif\"\n$pFiller.Range.InsertParagraphAfter()\n$pEndif = $d.Paragraphs.Item($count + 4)\n$pEndif.Range.InsertAfter(\"@endif\")\n"}
